# Fruta / hortaliza, semanal
#
# New weekly prices were added for "Ají" (Femacal de La Calera), which
# inserts two new rows right before the former row 675 and shifts every
# subsequent data row down by two (old row 675 -> new row 677, ...,
# old row 695 -> new row 697).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old row 675, pushing rows 675-695 down
# to 677-697 (Excel copies the row-675 formatting, incl. the date number
# format on column D, down onto the new rows automatically).
$ws.Rows("675:676").Insert()

# ---- New row 675 ----
$ws.Cells.Item(675, 1).Value = 3
$ws.Cells.Item(675, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(675, 3).Value = "Coquimbo"
$ws.Cells.Item(675, 4).Value2 = 45075
$ws.Cells.Item(675, 5).Value = 5
$ws.Cells.Item(675, 6).Value = 100112021
$ws.Cells.Item(675, 7).Value = "Ají"
$ws.Cells.Item(675, 8).Value = "Inferno"
$ws.Cells.Item(675, 9).Value = "Primera"
$ws.Cells.Item(675, 10).Value = 73
$ws.Cells.Item(675, 11).Value = 13500
$ws.Cells.Item(675, 12).Value = 14000
$ws.Cells.Item(675, 13).Value = 13760
$ws.Cells.Item(675, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(675, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(675, 16).Value = 1376
$ws.Cells.Item(675, 17).Value = 10
$ws.Cells.Item(675, 18).Value = "Hortaliza"

# ---- New row 676 ----
$ws.Cells.Item(676, 1).Value = 3
$ws.Cells.Item(676, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(676, 3).Value = "Coquimbo"
$ws.Cells.Item(676, 4).Value2 = 45075
$ws.Cells.Item(676, 5).Value = 5
$ws.Cells.Item(676, 6).Value = 100112021
$ws.Cells.Item(676, 7).Value = "Ají"
$ws.Cells.Item(676, 8).Value = "Inferno"
$ws.Cells.Item(676, 9).Value = "Primera"
$ws.Cells.Item(676, 10).Value = 73
$ws.Cells.Item(676, 11).Value = 23500
$ws.Cells.Item(676, 12).Value = 24000
$ws.Cells.Item(676, 13).Value = 23760
$ws.Cells.Item(676, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(676, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(676, 16).Value = 950
$ws.Cells.Item(676, 17).Value = 25
$ws.Cells.Item(676, 18).Value = "Hortaliza"
